# Remove the post row for "「行くな」ماتروحش" (originally row 738).
# Deleting the entire row shifts all subsequent rows (739-749) up by
# one, which also updates the sheet's used-range dimension automatically
# from A1:C749 to A1:C748.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(738).Delete()
